$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "38.734.09"
Set-TextValue $ws.Range("E2") "  +2.45%  "
Set-TextValue $ws.Range("D3") "2.084.11"
Set-TextValue $ws.Range("E3") "  +1.87%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "228.16"
Set-TextValue $ws.Range("E5") "  +0.24%  "
Set-TextValue $ws.Range("E6") "  +0.68%  "
Set-TextValue $ws.Range("D7") "60.24"
Set-TextValue $ws.Range("E7") "  +1.13%  "
Set-TextValue $ws.Range("E8") "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.383"
Set-TextValue $ws.Range("E9") "  +1.76%  "
Set-TextValue $ws.Range("D10") "0.0840"
Set-TextValue $ws.Range("E10") "  +0.30%  "
Set-TextValue $ws.Range("E11") "  -0.39%  "
Set-TextValue $ws.Range("D12") "2.393.55"
Set-TextValue $ws.Range("E12") "  +1.88%  "
Set-TextValue $ws.Range("D13") "14.96"
Set-TextValue $ws.Range("E13") "  +3.96%  "
Set-TextValue $ws.Range("D14") "21.84"
Set-TextValue $ws.Range("E14") "  +2.10%  "
Set-TextValue $ws.Range("E15") "  +4.16%  "
Set-TextValue $ws.Range("D16") "5.48"
Set-TextValue $ws.Range("E16") "  -0.38%  "
Set-TextValue $ws.Range("D17") "2.086.58"
Set-TextValue $ws.Range("E17") "  +2.11%  "
Set-TextValue $ws.Range("D18") "38.651.38"
Set-TextValue $ws.Range("E18") "  +2.34%  "
Set-TextValue $ws.Range("D19") "71.44"
Set-TextValue $ws.Range("E19") "  +2.83%  "
Set-TextValue $ws.Range("E20") "  +1.78%  "
Set-TextValue $ws.Range("D21") "0.0₃0839"
Set-TextValue $ws.Range("E21") "  +1.06%  "
Set-TextValue $ws.Range("D22") "227.19"
Set-TextValue $ws.Range("E22") "  +2.21%  "
Set-TextValue $ws.Range("D24") "2.43"
Set-TextValue $ws.Range("E24") "  +1.58%  "
Set-TextValue $ws.Range("E25") "  +1.99%  "
Set-TextValue $ws.Range("D26") "170.35"
Set-TextValue $ws.Range("E26") "  +1.05%  "
Set-TextValue $ws.Range("E27") "  +2.01%  "
Set-TextValue $ws.Range("E28") "  +8.72%  "
Set-TextValue $ws.Range("E29") "  +13.11%  "
Set-TextValue $ws.Range("E30") "  +2.02%  "
Set-TextValue $ws.Range("E31") "  +0.79%  "
Set-TextValue $ws.Range("E32") "  +5.94%  "
Set-TextValue $ws.Range("D33") "4.49"
Set-TextValue $ws.Range("E33") "  +2.73%  "
Set-TextValue $ws.Range("E34") "  +3.35%  "
Set-TextValue $ws.Range("D35") "0.0607"
Set-TextValue $ws.Range("E35") "  +0.97%  "
Set-TextValue $ws.Range("D36") "6.47"
Set-TextValue $ws.Range("E36") "  -0.25%  "
Set-TextValue $ws.Range("D37") "2.37"
Set-TextValue $ws.Range("E37") "  +1.35%  "
Set-TextValue $ws.Range("D38") "3.54"
Set-TextValue $ws.Range("E38") "  +2.04%  "
Set-TextValue $ws.Range("E39") "  -0.13%  "
Set-TextValue $ws.Range("D40") "17.91"
Set-TextValue $ws.Range("E40") "  -2.39%  "
Set-TextValue $ws.Range("D41") "0.0225"
Set-TextValue $ws.Range("E41") "  +4.66%  "
Set-TextValue $ws.Range("D42") "1.537.82"
Set-TextValue $ws.Range("E42") "  +0.92%  "
Set-TextValue $ws.Range("D43") "100.46"
Set-TextValue $ws.Range("E43") "  +3.00%  "
Set-TextValue $ws.Range("E44") "  -0.65%  "
Set-TextValue $ws.Range("E45") "  +3.43%  "
Set-TextValue $ws.Range("E46") "  +8.17%  "
Set-TextValue $ws.Range("E47") "  +1.41%  "
Set-TextValue $ws.Range("E48") "  -2.52%  "
Set-TextValue $ws.Range("E49") "  +2.51%  "
Set-TextValue $ws.Range("E50") "  +0.54%  "
Set-TextValue $ws.Range("D51") "2.282.79"
Set-TextValue $ws.Range("E51") "  +2.03%  "
